$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting from row 4 onto the new row 5 so styles (borders / wrap /
# hyperlink font) line up with the existing data rows before the values are
# written.
$ws.Range("B4:M4").Copy()
$ws.Range("B5:M5").PasteSpecial(-4122)

# Populate the new FAQ row.
$ws.Range("A5").Value = "No Spaces Warning"
$ws.Range("B5").Value = "What is Q and A Bot"
$ws.Range("C5").Value = "What is QnaBot"
$ws.Range("D5").Value = "The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Range("E5").Value = "The Q and A Bot uses [Amazon Lex](https://aws.amazon.com/lex/) and [Alexa](https://developer.amazon.com/en-US/alexa) to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Range("F5").Value = "<speak>The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer</speak>"
$ws.Range("G5").Value = "Alexa"
$ws.Range("H5").Value = "Alexa"
$ws.Range("I5").Value = "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg"
$ws.Range("J5").Value = "Tell me about the Alexa Show."
$ws.Range("K5").Value = "The Echo Show"
$ws.Range("L5").Value = "Tell me about the Echo Dot"
$ws.Range("M5").Value = "The Echo Dot"

# Row 5 has wrapped long text, so it gets a tall auto row height like rows 2-4.
$ws.Rows.Item(5).RowHeight = 153

# Hyperlink the new image-URL cell like the rows above it.
$ws.Hyperlinks.Add($ws.Range("I5"), "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg")

# Scroll the sheet so row 2 is pinned at top and the new row is selected,
# matching the saved view state.
$ws.Range("B5:M5").Select()
$ws.Application.ActiveWindow.ScrollRow = 2

$excel.ActiveWindow.WindowState = -4137
$excel.Width = 30240
$excel.Height = 18880
